$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape run
# Plain decimal "Price" values must be forced to text so Excel
# does not reinterpret them as numbers (avoids float rounding and
# preserves the original inline-string/general-format cell shape).
$ws.Range("D2").Value = '67.822.40'
$ws.Range("E2").Value = '  +3.37%  '
$ws.Range("D3").Value = '3.292.55'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  +0.01%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '578.41'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.27%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '179.37'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -1.44%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.11%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.587'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +3.13%  '
$ws.Range("D9").Value = '3.284.40'
$ws.Range("E9").Value = '  +0.21%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.177'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("E11").Value = '  +1.11%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '45.83'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("E13").Value = '  +3.33%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '692.52'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +12.81%  '
$ws.Range("D15").Value = '3.815.94'
$ws.Range("E15").Value = '  +0.42%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '8.40'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").Value = '67.894.51'
$ws.Range("E17").Value = '  +3.31%  '
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("D19").Value = '3.293.76'
$ws.Range("E19").Value = '  +0.72%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '17.47'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -1.18%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '10.82'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -0.37%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.896'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +1.37%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '17.09'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -4.92%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '5.21'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +5.74%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '98.59'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.46%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '3.95'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.19%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '2.75'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +1.86%  '
$ws.Range("E28").Value = '  -0.30%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '33.18'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +7.98%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '8.49'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +1.77%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '6.79'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +5.30%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '582.99'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +7.28%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '10.90'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.84%  '
$ws.Range("D34").Value = '3.867.02'
$ws.Range("E34").Value = '  +2.00%  '
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("E36").Value = '  +0.14%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '3.40'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -8.30%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '55.37'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -1.03%  '
$ws.Range("E39").Value = '  +1.82%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '3.22'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +3.01%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '2.63'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +2.57%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '32.30'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '3.39'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("D44").Value = '0.0₃0680'
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("E46").Value = '  +2.19%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.129'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +2.36%  '
$ws.Range("E48").Value = '  +10.40%  '
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("E50").Value = '  +1.32%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '129.27'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +1.16%  '
